$wb = $excel.ActiveWorkbook

# ==================== Sheet 1: Overview ====================
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 6 (shifts rows 6-8 down to 7-9,
# copies formatting from the row above automatically)
$ws.Rows.Item(6).Insert()

# Populate the new row with the handed-off files data
$ws.Range("A6").Value = '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md'
$ws.Range("B6").Value = 'Ready for handoff'
$ws.Range("C6").Value = 'Ready for handoff'

# Rebuild the hyperlinks collection in correct row order
# (Insert() does not shift existing hyperlink row bindings)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/39b27abd8ea5205bfeaefac4e8f381dd7ecfd1c9/e2e/8170b961-4bbf-4106-ad72-790e4f0e3c4e.md', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/13e9c5860886f944ea821f40324f3d75344f9c2f/e2e/d7959a83-2677-4c73-8290-262faedc1f7d.md', "", "", 'd7959a83-2677-4c73-8290-262faedc1f7d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/7010838e01a097bf9537c8a9cafb9a4381564c04/e2e/ddae8de3-6496-4ed0-9fcb-eef710189763.md', "", "", 'ddae8de3-6496-4ed0-9fcb-eef710189763.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), 'https://github.com/OpenLocalizationTest/oltest/blob/7010838e01a097bf9537c8a9cafb9a4381564c04/e2e/e390d7bc-7fae-43e5-83b5-43910950d9d8.md', "", "", 'e390d7bc-7fae-43e5-83b5-43910950d9d8.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), 'https://github.com/OpenLocalizationTest/oltest/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/e2e/13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md', "", "", '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), 'https://github.com/OpenLocalizationTest/oltest/blob/55a38757d76db317366e741e705226dd0a0e8e2c/e2e/4f748b75-aebf-4922-8893-6d729490a2c4.md', "", "", '4f748b75-aebf-4922-8893-6d729490a2c4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), 'https://github.com/OpenLocalizationTest/oltest/blob/88470dcd45a9605e17d349330815d11dc9253709/e2e/c8414e6c-ed0c-400d-8a83-8089a6a4569f.md', "", "", 'c8414e6c-ed0c-400d-8a83-8089a6a4569f.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), 'https://github.com/OpenLocalizationTest/oltest/blob/88470dcd45a9605e17d349330815d11dc9253709/.localization-config', "", "", '.localization-config') | Out-Null

# ==================== Sheet 2: zh-cn ====================
$ws = $wb.Worksheets.Item(2)

# Insert a new row at position 6 (shifts rows 6-8 down to 7-9,
# copies formatting from the row above automatically)
$ws.Rows.Item(6).Insert()

# Populate the new row with the handed-off files data
$ws.Range("A6").Value = '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md'
$ws.Range("B6").Value = 'Ready for handoff'
$ws.Range("C6").Value = '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.f554c1a12e8d435602f3963fb78e5c4babf0b95e.zh-cn.xlf'
$ws.Range("D6").Value = '2016-02-18 09:36:20'
$ws.Range("G6").Value = '0001-01-01 00:00:00'
$ws.Range("H6").Value = 'Include'

# Rebuild the hyperlinks collection in correct row order
# (Insert() does not shift existing hyperlink row bindings)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/39b27abd8ea5205bfeaefac4e8f381dd7ecfd1c9/e2e/8170b961-4bbf-4106-ad72-790e4f0e3c4e.md', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1f6657de62b6b6bd7ffbdac3e70c89317643988/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.zh-cn.xlf', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/feff45639d5a1a10ae6df4a7d72a5d283e5315ab/e2e/8170b961-4bbf-4106-ad72-790e4f0e3c4e.md', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/fc5337cb1c5d7893a64ae5eec11dbec28d9d2806/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.zh-cn.xlf', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/13e9c5860886f944ea821f40324f3d75344f9c2f/e2e/d7959a83-2677-4c73-8290-262faedc1f7d.md', "", "", 'd7959a83-2677-4c73-8290-262faedc1f7d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e61078211c5b355eec014f16662bc8ee2b85d17e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d7959a83-2677-4c73-8290-262faedc1f7d.9a38d13c4ece2e3926f3f99a235bf2b62096fa98.zh-cn.xlf', "", "", 'd7959a83-2677-4c73-8290-262faedc1f7d.9a38d13c4ece2e3926f3f99a235bf2b62096fa98.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/7010838e01a097bf9537c8a9cafb9a4381564c04/e2e/ddae8de3-6496-4ed0-9fcb-eef710189763.md', "", "", 'ddae8de3-6496-4ed0-9fcb-eef710189763.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/04908536571ddef24653d977338b769779a9472f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ddae8de3-6496-4ed0-9fcb-eef710189763.b8a12fe613db63046f51e7bb82776e262d50f1cc.zh-cn.xlf', "", "", 'ddae8de3-6496-4ed0-9fcb-eef710189763.b8a12fe613db63046f51e7bb82776e262d50f1cc.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), 'https://github.com/OpenLocalizationTest/oltest/blob/7010838e01a097bf9537c8a9cafb9a4381564c04/e2e/e390d7bc-7fae-43e5-83b5-43910950d9d8.md', "", "", 'e390d7bc-7fae-43e5-83b5-43910950d9d8.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/04908536571ddef24653d977338b769779a9472f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e390d7bc-7fae-43e5-83b5-43910950d9d8.895648431ff115c06344d9005327b8f6de8d963d.zh-cn.xlf', "", "", 'e390d7bc-7fae-43e5-83b5-43910950d9d8.895648431ff115c06344d9005327b8f6de8d963d.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), 'https://github.com/OpenLocalizationTest/oltest/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/e2e/13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md', "", "", '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.f554c1a12e8d435602f3963fb78e5c4babf0b95e.zh-cn.xlf', "", "", '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.f554c1a12e8d435602f3963fb78e5c4babf0b95e.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), 'https://github.com/OpenLocalizationTest/oltest/blob/55a38757d76db317366e741e705226dd0a0e8e2c/e2e/4f748b75-aebf-4922-8893-6d729490a2c4.md', "", "", '4f748b75-aebf-4922-8893-6d729490a2c4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/142bb5a56547cd028775252a97f348348d51decd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4f748b75-aebf-4922-8893-6d729490a2c4.cb43d0086ff54aab8af6c518058b4d561f411e27.zh-cn.xlf', "", "", '4f748b75-aebf-4922-8893-6d729490a2c4.cb43d0086ff54aab8af6c518058b4d561f411e27.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), 'https://github.com/OpenLocalizationTest/oltest/blob/88470dcd45a9605e17d349330815d11dc9253709/e2e/c8414e6c-ed0c-400d-8a83-8089a6a4569f.md', "", "", 'c8414e6c-ed0c-400d-8a83-8089a6a4569f.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7d5ed50ad9d8418dea0fe518fa78b99700b626f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c8414e6c-ed0c-400d-8a83-8089a6a4569f.6cd3cd2ec43daddffcd4a76ad36e6adfef93894a.zh-cn.xlf', "", "", 'c8414e6c-ed0c-400d-8a83-8089a6a4569f.6cd3cd2ec43daddffcd4a76ad36e6adfef93894a.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), 'https://github.com/OpenLocalizationTest/oltest/blob/88470dcd45a9605e17d349330815d11dc9253709/.localization-config', "", "", '.localization-config') | Out-Null

# ==================== Sheet 3: de-de ====================
$ws = $wb.Worksheets.Item(3)

# Insert a new row at position 6 (shifts rows 6-8 down to 7-9,
# copies formatting from the row above automatically)
$ws.Rows.Item(6).Insert()

# Populate the new row with the handed-off files data
$ws.Range("A6").Value = '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md'
$ws.Range("B6").Value = 'Ready for handoff'
$ws.Range("C6").Value = '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.f554c1a12e8d435602f3963fb78e5c4babf0b95e.de-de.xlf'
$ws.Range("D6").Value = '2016-02-18 09:36:31'
$ws.Range("G6").Value = '0001-01-01 00:00:00'
$ws.Range("H6").Value = 'Include'

# Rebuild the hyperlinks collection in correct row order
# (Insert() does not shift existing hyperlink row bindings)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/39b27abd8ea5205bfeaefac4e8f381dd7ecfd1c9/e2e/8170b961-4bbf-4106-ad72-790e4f0e3c4e.md', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4f674a30a4130844d06338f7b6e06270360b7ee8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.de-de.xlf', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/93fa94d4f47b4ebe3fcf0ea49c1c851b0be4d472/e2e/8170b961-4bbf-4106-ad72-790e4f0e3c4e.md', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/e13ddcead133549a0e8dd0fd3dc802ef0674d324/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.de-de.xlf', "", "", '8170b961-4bbf-4106-ad72-790e4f0e3c4e.be393b1af074b9cc113a7e2ed2c9604eb7a6de0e.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/13e9c5860886f944ea821f40324f3d75344f9c2f/e2e/d7959a83-2677-4c73-8290-262faedc1f7d.md', "", "", 'd7959a83-2677-4c73-8290-262faedc1f7d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bad00aaf0bde640cac6bde45a045513d32d4cf9c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/d7959a83-2677-4c73-8290-262faedc1f7d.9a38d13c4ece2e3926f3f99a235bf2b62096fa98.de-de.xlf', "", "", 'd7959a83-2677-4c73-8290-262faedc1f7d.9a38d13c4ece2e3926f3f99a235bf2b62096fa98.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/7010838e01a097bf9537c8a9cafb9a4381564c04/e2e/ddae8de3-6496-4ed0-9fcb-eef710189763.md', "", "", 'ddae8de3-6496-4ed0-9fcb-eef710189763.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e4a8b69587296e54c0f494b4018a1fe8a404ecc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ddae8de3-6496-4ed0-9fcb-eef710189763.b8a12fe613db63046f51e7bb82776e262d50f1cc.de-de.xlf', "", "", 'ddae8de3-6496-4ed0-9fcb-eef710189763.b8a12fe613db63046f51e7bb82776e262d50f1cc.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), 'https://github.com/OpenLocalizationTest/oltest/blob/7010838e01a097bf9537c8a9cafb9a4381564c04/e2e/e390d7bc-7fae-43e5-83b5-43910950d9d8.md', "", "", 'e390d7bc-7fae-43e5-83b5-43910950d9d8.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e4a8b69587296e54c0f494b4018a1fe8a404ecc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e390d7bc-7fae-43e5-83b5-43910950d9d8.895648431ff115c06344d9005327b8f6de8d963d.de-de.xlf', "", "", 'e390d7bc-7fae-43e5-83b5-43910950d9d8.895648431ff115c06344d9005327b8f6de8d963d.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), 'https://github.com/OpenLocalizationTest/oltest/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/e2e/13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md', "", "", '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.f554c1a12e8d435602f3963fb78e5c4babf0b95e.de-de.xlf', "", "", '13384a12-7135-4fc0-a9ff-09e0a8ea0c6a.f554c1a12e8d435602f3963fb78e5c4babf0b95e.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), 'https://github.com/OpenLocalizationTest/oltest/blob/55a38757d76db317366e741e705226dd0a0e8e2c/e2e/4f748b75-aebf-4922-8893-6d729490a2c4.md', "", "", '4f748b75-aebf-4922-8893-6d729490a2c4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e350e68a1bc54ef9912571a0b2bdda9fbc7cb553/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4f748b75-aebf-4922-8893-6d729490a2c4.cb43d0086ff54aab8af6c518058b4d561f411e27.de-de.xlf', "", "", '4f748b75-aebf-4922-8893-6d729490a2c4.cb43d0086ff54aab8af6c518058b4d561f411e27.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), 'https://github.com/OpenLocalizationTest/oltest/blob/88470dcd45a9605e17d349330815d11dc9253709/e2e/c8414e6c-ed0c-400d-8a83-8089a6a4569f.md', "", "", 'c8414e6c-ed0c-400d-8a83-8089a6a4569f.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4a7d38d906faafcb42244fcf3da2bfc87d81740/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c8414e6c-ed0c-400d-8a83-8089a6a4569f.6cd3cd2ec43daddffcd4a76ad36e6adfef93894a.de-de.xlf', "", "", 'c8414e6c-ed0c-400d-8a83-8089a6a4569f.6cd3cd2ec43daddffcd4a76ad36e6adfef93894a.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), 'https://github.com/OpenLocalizationTest/oltest/blob/88470dcd45a9605e17d349330815d11dc9253709/.localization-config', "", "", '.localization-config') | Out-Null
